$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto market data rows (price + 1h volume change columns),
# and fix the Hedera/Filecoin row ordering, per the latest pull.

$ws.Range("D2").Value = '58.037.73'
$ws.Range("E2").Value = '  +1.42%  '

$ws.Range("D3").Value = '3.134.51'
$ws.Range("E3").Value = '  +1.35%  '

$ws.Range("E4").Value = '  -0.03%  '

$fmtD5 = $ws.Range("D5").NumberFormat
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.32'
$ws.Range("D5").NumberFormat = $fmtD5
$ws.Range("E5").Value = '  +2.28%  '

$fmtD6 = $ws.Range("D6").NumberFormat
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.86'
$ws.Range("D6").NumberFormat = $fmtD6
$ws.Range("E6").Value = '  +1.85%  '

$fmtD7 = $ws.Range("D7").NumberFormat
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").NumberFormat = $fmtD7
$ws.Range("E7").Value = '  +0.07%  '

$fmtD8 = $ws.Range("D8").NumberFormat
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.510'
$ws.Range("D8").NumberFormat = $fmtD8
$ws.Range("E8").Value = '  +11.38%  '

$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("E10").Value = '  +2.11%  '

$fmtD11 = $ws.Range("D11").NumberFormat
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.419'
$ws.Range("D11").NumberFormat = $fmtD11
$ws.Range("E11").Value = '  +4.67%  '

$ws.Range("E12").Value = '  +3.47%  '

$ws.Range("D13").Value = '3.675.31'
$ws.Range("E13").Value = '  +1.32%  '

$ws.Range("E14").Value = '  +1.07%  '

$fmtD15 = $ws.Range("D15").NumberFormat
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000170'
$ws.Range("D15").NumberFormat = $fmtD15
$ws.Range("E15").Value = '  +5.13%  '

$ws.Range("D16").Value = '58.094.16'

$fmtD17 = $ws.Range("D17").NumberFormat
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.24'
$ws.Range("D17").NumberFormat = $fmtD17
$ws.Range("E17").Value = '  +5.81%  '

$ws.Range("D18").Value = '3.142.65'
$ws.Range("E18").Value = '  +1.61%  '

$fmtD19 = $ws.Range("D19").NumberFormat
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.95'
$ws.Range("D19").NumberFormat = $fmtD19
$ws.Range("E19").Value = '  +3.79%  '

$fmtD20 = $ws.Range("D20").NumberFormat
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.19'
$ws.Range("D20").NumberFormat = $fmtD20
$ws.Range("E20").Value = '  +4.09%  '

$fmtD21 = $ws.Range("D21").NumberFormat
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.60'
$ws.Range("D21").NumberFormat = $fmtD21
$ws.Range("E21").Value = '  +7.18%  '

$fmtD22 = $ws.Range("D22").NumberFormat
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").NumberFormat = $fmtD22
$ws.Range("E22").Value = '  +0.02%  '

$fmtD23 = $ws.Range("D23").NumberFormat
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.74'
$ws.Range("D23").NumberFormat = $fmtD23
$ws.Range("E23").Value = '  -0.79%  '

$fmtD24 = $ws.Range("D24").NumberFormat
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.10'
$ws.Range("D24").NumberFormat = $fmtD24
$ws.Range("E24").Value = '  +2.10%  '

$fmtD25 = $ws.Range("D25").NumberFormat
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.514'
$ws.Range("D25").NumberFormat = $fmtD25
$ws.Range("E25").Value = '  +2.97%  '

$ws.Range("E26").Value = '  +0.31%  '

$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D28").Value = '0.0₃0884'
$ws.Range("E28").Value = '  +1.94%  '

$fmtD29 = $ws.Range("D29").NumberFormat
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.85'
$ws.Range("D29").NumberFormat = $fmtD29
$ws.Range("E29").Value = '  +7.95%  '

$fmtD30 = $ws.Range("D30").NumberFormat
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.19'
$ws.Range("D30").NumberFormat = $fmtD30
$ws.Range("E30").Value = '  +6.04%  '

$ws.Range("E31").Value = '  +0.90%  '

$fmtD32 = $ws.Range("D32").NumberFormat
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.82'
$ws.Range("D32").NumberFormat = $fmtD32
$ws.Range("E32").Value = '  +4.34%  '

$fmtD33 = $ws.Range("D33").NumberFormat
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.19'
$ws.Range("D33").NumberFormat = $fmtD33
$ws.Range("E33").Value = '  +6.55%  '

$ws.Range("E34").Value = '  +3.18%  '

$fmtD35 = $ws.Range("D35").NumberFormat
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '161.27'
$ws.Range("D35").NumberFormat = $fmtD35
$ws.Range("E35").Value = '  +1.54%  '

$ws.Range("E36").Value = '  +3.47%  '

$ws.Range("E37").Value = '  +8.68%  '

$fmtD38 = $ws.Range("D38").NumberFormat
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.54'
$ws.Range("D38").NumberFormat = $fmtD38
$ws.Range("E38").Value = '  -0.17%  '

$ws.Range("E39").Value = '  +6.34%  '

$ws.Range("D40").Value = '2.624.85'
$ws.Range("E40").Value = '  +9.40%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$fmtD41 = $ws.Range("D41").NumberFormat
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.22'
$ws.Range("D41").NumberFormat = $fmtD41
$ws.Range("E41").Value = '  +4.35%  '

$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$fmtD42 = $ws.Range("D42").NumberFormat
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0674'
$ws.Range("D42").NumberFormat = $fmtD42
$ws.Range("E42").Value = '  +2.37%  '

$fmtD43 = $ws.Range("D43").NumberFormat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.91'
$ws.Range("D43").NumberFormat = $fmtD43
$ws.Range("E43").Value = '  +6.05%  '

$ws.Range("E44").Value = '  +0.67%  '

$ws.Range("E45").Value = '  +3.19%  '

$fmtD46 = $ws.Range("D46").NumberFormat
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").NumberFormat = $fmtD46
$ws.Range("E46").Value = '  -0.05%  '

$ws.Range("E47").Value = '  +4.30%  '

$ws.Range("E48").Value = '  +2.99%  '

$fmtD49 = $ws.Range("D49").NumberFormat
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0994'
$ws.Range("D49").NumberFormat = $fmtD49
$ws.Range("E49").Value = '  +10.09%  '

$fmtD50 = $ws.Range("D50").NumberFormat
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.27'
$ws.Range("D50").NumberFormat = $fmtD50
$ws.Range("E50").Value = '  +2.90%  '

$fmtD51 = $ws.Range("D51").NumberFormat
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.750'
$ws.Range("D51").NumberFormat = $fmtD51
$ws.Range("E51").Value = '  -1.98%  '
